$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank row that was row 9 (the row below the "Cenário Normal" header block)
$ws.Rows(9).Delete()

# Update text content
$ws.Range("C2").Value = "Apresenta informação extra"
$ws.Range("C3").Value = "Utilizador autenticado"
$ws.Range("C4").Value = "A confecionar passo"

# Update font used across the table: Calibri 14 -> Arial 11
$ws.Range("B2:D16").Font.Name = "Arial"
$ws.Range("B2:D16").Font.Size = 11

# Selection / view cosmetics
$ws.Range("C23").Select()
